# Regenerate merged AHB files
# 1) Rename header labels: "<Name>_old" -> "<Name>_FV2310" and "<Name>_new" -> "<Name>_FV2404"
# 2) Turn the data range into an Excel Table (ListObject) with an AutoFilter
# 3) Freeze the header row (pane split after row 1)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the header row labels -------------------------------------
$lastCol = 21   # A..U
for ($col = 1; $col -le $lastCol; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $val = $cell.Value2
    if ($val -ne $null) {
        if ($val -like "*_old") {
            $cell.Value = ($val -replace "_old$", "_FV2310")
        } elseif ($val -like "*_new") {
            $cell.Value = ($val -replace "_new$", "_FV2404")
        }
    }
}

# --- 2. Create the table over A1:U89 --------------------------------------
$tableRange = $ws.Range("A1:U89")
$lo = $ws.ListObjects.Add(1, $tableRange, $true, 1, "Table1")
$lo.Name = "Table1"

# --- 3. Freeze panes at row 1 ----------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
